# Apply updated cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.402.44'
$ws.Range('E2').Value = '  +0.87%  '
# Row 3
$ws.Range('D3').Value = '1.997.78'
$ws.Range('E3').Value = '  +3.83%  '
# Row 4
$ws.Range('E4').Value = '  -0.02%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.51'
$ws.Range('E5').Value = '  +0.76%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.06%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5108'
$ws.Range('E7').Value = '  +0.87%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4130'
$ws.Range('E8').Value = '  +2.23%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08714'
$ws.Range('E9').Value = '  +5.25%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.135'
$ws.Range('E10').Value = '  +1.95%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.74'
$ws.Range('E11').Value = '  +1.73%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.75'
$ws.Range('E12').Value = '  +3.59%  '
# Row 13
$ws.Range('D13').Value = '1.998.44'
$ws.Range('E13').Value = '  +4.84%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.537'
$ws.Range('E14').Value = '  +1.72%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.437'
$ws.Range('E15').Value = '  +1.30%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9999'
$ws.Range('E16').Value = '  -0.09%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.07'
$ws.Range('E17').Value = '  +1.34%  '
# Row 18
$ws.Range('E18').Value = '  +1.26%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06504'
$ws.Range('E19').Value = '  +0.26%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.93'
$ws.Range('E20').Value = '  +1.49%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9990'
$ws.Range('E21').Value = '  -0.16%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.172'
$ws.Range('E22').Value = '  +3.03%  '
# Row 23
$ws.Range('D23').Value = '30.458.58'
$ws.Range('E23').Value = '  +0.92%  '
# Row 24
$ws.Range('E24').Value = '  +4.58%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.207'
$ws.Range('E25').Value = '  +0.56%  '
# Row 26
$ws.Range('D26').Value = '2.232.71'
$ws.Range('E26').Value = '  +5.04%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.45'
$ws.Range('E27').Value = '  +0.91%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.17'
$ws.Range('E28').Value = '  +1.35%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.420'
$ws.Range('E29').Value = '  +1.84%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.51'
$ws.Range('E30').Value = '  +1.49%  '
# Row 31
$ws.Range('E31').Value = '  +0.96%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1051'
$ws.Range('E32').Value = '  +0.54%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.071'
$ws.Range('E33').Value = '  +0.83%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.834'
$ws.Range('E34').Value = '  +1.07%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.332'
$ws.Range('E35').Value = '  +11.76%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02510'
$ws.Range('E36').Value = '  +2.32%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06591'
$ws.Range('E37').Value = '  +1.89%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.365'
$ws.Range('E38').Value = '  -1.54%  '
# Row 39
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2198'
$ws.Range('E39').Value = '  +1.53%  '
# Row 40
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.18'
$ws.Range('E40').Value = '  +6.77%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.023'
$ws.Range('E41').Value = '  +1.81%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6627'
$ws.Range('E42').Value = '  +3.37%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.228'
$ws.Range('E43').Value = '  +0.85%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.66'
$ws.Range('E44').Value = '  +2.72%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6162'
$ws.Range('E45').Value = '  +2.50%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.205'
$ws.Range('E46').Value = '  +1.70%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.662'
$ws.Range('E47').Value = '  +0.34%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.259'
$ws.Range('E48').Value = '  +3.32%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.56'
$ws.Range('E49').Value = '  +0.54%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '80.36'
$ws.Range('E50').Value = '  +1.36%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06892'
$ws.Range('E51').Value = '  +1.13%  '
